$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41 (G41=5478)
$ws.Range("H41").Value = 723.2
$ws.Range("I41").Value = 107.5
$ws.Range("K41").Value = 107.5
$ws.Range("M41").Value = 332.5

# Row 43 (G43=5472)
$ws.Range("H43").Value = 1712.8572
$ws.Range("J43").Value = 2000
$ws.Range("L43").Value = 2000
$ws.Range("N43").Value = -2138

# Row 51 (G51=5486)
$ws.Range("H51").Value = 3697
$ws.Range("J51").Value = 3700
$ws.Range("L51").Value = 3700
$ws.Range("N51").Value = -4668

# Row 55 (G55=5517)
$ws.Range("H55").Value = 329.42856
$ws.Range("I55").Value = 655.5
$ws.Range("J55").Value = 199
$ws.Range("K55").Value = 655.5
$ws.Range("L55").Value = 199
$ws.Range("M55").Value = -441.5
$ws.Range("N55").Value = -627

# Row 132 (G132=44049)
$ws.Range("H132").Value = 1202.375
$ws.Range("I132").Value = 1202.375
$ws.Range("K132").Value = 3607.125
$ws.Range("M132").Value = -1077.125

$ws = $wb.Worksheets.Item("ARM")
# Row 5 (G5=5091)
$ws.Range("H5").Value = 1608.5714
$ws.Range("I5").Value = 1608.5714
$ws.Range("K5").Value = 1608.5714
$ws.Range("M5").Value = -1496.5714

# Row 8 (G8=3011)
$ws.Range("H8").Value = 406560
$ws.Range("I8").Value = 669933.3
$ws.Range("J8").Value = 11500
$ws.Range("K8").Value = 669933.3
$ws.Range("L8").Value = 11500
$ws.Range("M8").Value = -669789.3
$ws.Range("N8").Value = -11788

# Row 32 (G32=44147)
$ws.Range("H32").Value = 6942.1333
$ws.Range("I32").Value = 7080.857
$ws.Range("K32").Value = 7080.857
$ws.Range("M32").Value = -6793.857

# Row 41 (G41=2501)
$ws.Range("H41").Value = 1997.6666
$ws.Range("I41").Value = 1997.6666
$ws.Range("K41").Value = 1997.6666
$ws.Range("M41").Value = -1583.6666

# Row 63 (G63=12528)
$ws.Range("H63").Value = 25981
$ws.Range("I63").Value = 25981
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 25981
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -25295
$ws.Range("N63").ClearContents()

# Row 66 (G66=12528)
$ws.Range("H66").Value = 25981
$ws.Range("I66").Value = 25981
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 129905
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -126473
$ws.Range("N66").ClearContents()

# Row 88 (G88=12530)
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# Row 91 (G91=12530)
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 4 (G4=5091)
$ws.Range("H4").Value = 1608.5714
$ws.Range("I4").Value = 1608.5714
$ws.Range("K4").Value = 1608.5714
$ws.Range("M4").Value = -1493.5714

# Row 22 (G22=5092)
$ws.Range("H22").Value = 840.8889
$ws.Range("I22").Value = 547.8333
$ws.Range("J22").Value = 1427
$ws.Range("K22").Value = 547.8333
$ws.Range("L22").Value = 1427
$ws.Range("M22").Value = -374.8333
$ws.Range("N22").Value = -1773

# Row 130 (G130=34682)
$ws.Range("H130").Value = 99998.55499999999
$ws.Range("J130").Value = 99998.55499999999
$ws.Range("L130").Value = 99998.55499999999
$ws.Range("N130").Value = -110038.555

$ws = $wb.Worksheets.Item("CRP")
# Row 12 (G12=1604)
$ws.Range("H12").Value = 5355.5
$ws.Range("J12").Value = 706
$ws.Range("L12").Value = 706
$ws.Range("N12").Value = -1046

# Row 31 (G31=44023)
$ws.Range("H31").Value = 1639.9166
$ws.Range("I31").Value = 1390.2
$ws.Range("J31").Value = 1818.2858
$ws.Range("K31").Value = 1390.2
$ws.Range("L31").Value = 1818.2858
$ws.Range("M31").Value = -1095.2
$ws.Range("N31").Value = -2408.2858

# Row 34 (G34=44023)
$ws.Range("H34").Value = 1639.9166
$ws.Range("I34").Value = 1390.2
$ws.Range("J34").Value = 1818.2858
$ws.Range("K34").Value = 1390.2
$ws.Range("L34").Value = 1818.2858
$ws.Range("M34").Value = -1188.2
$ws.Range("N34").Value = -2222.2858

# Row 134 (G134=44020)
$ws.Range("H134").Value = 2098.0667
$ws.Range("I134").Value = 2479.2727
$ws.Range("K134").Value = 7437.8181
$ws.Range("M134").Value = -4902.8181

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (G4=4650)
$ws.Range("H4").Value = 15715159
$ws.Range("I4").Value = 15715159
$ws.Range("K4").Value = 47145477
$ws.Range("M4").Value = -47145365

# Row 14 (G14=12886)
$ws.Range("H14").Value = 900.6667
$ws.Range("I14").Value = 900.6667
$ws.Range("K14").Value = 2702.0001
$ws.Range("M14").Value = -2529.0001

# Row 32 (G32=4731)
$ws.Range("H32").Value = 1300
$ws.Range("I32").Value = 1300
$ws.Range("K32").Value = 3900
$ws.Range("M32").Value = -3617

# Row 46 (G46=4701)
$ws.Range("H46").Value = 449.57144
$ws.Range("I46").Value = 45
$ws.Range("J46").Value = 517
$ws.Range("K46").Value = 135
$ws.Range("L46").Value = 1551
$ws.Range("M46").Value = -44
$ws.Range("N46").Value = -1733

# Row 122 (G122=36078)
$ws.Range("H122").Value = 899.6667
$ws.Range("J122").Value = 899.6667
$ws.Range("L122").Value = 8097.0003
$ws.Range("N122").Value = -12997.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 10 (G10=4306)
$ws.Range("H10").Value = 350.5
$ws.Range("I10").Value = 350.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 350.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -181.5
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 5 (G5=3790)
$ws.Range("H5").Value = 251250
$ws.Range("I5").Value = 2500
$ws.Range("J5").Value = 500000
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 500000
$ws.Range("M5").Value = -2387
$ws.Range("N5").Value = -500226

# Row 18 (G18=3772)
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 50
$ws.Range("K18").Value = 50
$ws.Range("M18").Value = 122

# Row 19 (G19=2229)
$ws.Range("H19").Value = 533.3333
$ws.Range("I19").Value = 533.3333
$ws.Range("K19").Value = 533.3333
$ws.Range("M19").Value = -363.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 17 (G17=3539)
$ws.Range("H17").Value = 404
$ws.Range("I17").Value = 404
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 404
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -232
$ws.Range("N17").ClearContents()

# Row 119 (G119=26289)
$ws.Range("H119").Value = 63750
$ws.Range("J119").Value = 63750
$ws.Range("L119").Value = 63750
$ws.Range("N119").Value = -73426

# Row 122 (G122=36208)
$ws.Range("H122").Value = 962.25
$ws.Range("J122").Value = 900
$ws.Range("L122").Value = 2700
$ws.Range("N122").Value = -7600

# Row 132 (G132=44029)
$ws.Range("H132").Value = 6179.231
$ws.Range("I132").Value = 6376.4287
$ws.Range("K132").Value = 19129.2861
$ws.Range("M132").Value = -16599.2861

# Row 136 (G136=44031)
$ws.Range("H136").Value = 3673.5652
$ws.Range("I136").Value = 3541.875
$ws.Range("K136").Value = 10625.625
$ws.Range("M136").Value = -8075.625

